$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header in B1 from "2015 a 2022" to "Valor"
$ws.Range("B1").Value = "Valor"

# Reflect the final selection state seen in the saved file (single cell B2)
$ws.Range("B2").Select()
